# Update the date heading and the division expressions in the table.
$d = $word.ActiveDocument

$replacements = @(
    @{ Old = "2025-09-29 Monday"; New = "2025-09-30 Tuesday" },
    @{ Old = "717÷2="; New = "463÷7=" },
    @{ Old = "431÷5="; New = "708÷5=" },
    @{ Old = "566÷7="; New = "414÷7=" },
    @{ Old = "914÷7="; New = "173÷9=" },
    @{ Old = "264÷7="; New = "324÷8=" },
    @{ Old = "874÷7="; New = "429÷6=" },
    @{ Old = "212÷6="; New = "526÷5=" },
    @{ Old = "387÷9="; New = "455÷5=" },
    @{ Old = "648÷6="; New = "437÷7=" },
    @{ Old = "402÷9="; New = "489÷2=" },
    @{ Old = "651÷2="; New = "883÷2=" },
    @{ Old = "995÷4="; New = "340÷5=" },
    @{ Old = "781÷6="; New = "580÷3=" },
    @{ Old = "698÷5="; New = "310÷2=" },
    @{ Old = "157÷4="; New = "647÷8=" },
    @{ Old = "695÷9="; New = "519÷5=" },
    @{ Old = "870÷7="; New = "123÷4=" },
    @{ Old = "461÷2="; New = "499÷9=" },
    @{ Old = "857÷2="; New = "548÷2=" },
    @{ Old = "639÷6="; New = "833÷2=" },
    @{ Old = "694÷7="; New = "787÷9=" },
    @{ Old = "279÷5="; New = "894÷2=" },
    @{ Old = "540÷5="; New = "173÷4=" },
    @{ Old = "282÷8="; New = "519÷4=" },
    @{ Old = "830÷7="; New = "761÷5=" }
)

foreach ($pair in $replacements) {
    $d.Content.Find.Execute($pair.Old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $pair.New, 2)
}
